$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planning")
$ws.Activate()

# Rows 15-25 in column B (Activiteit) get new / renamed text, column D (Uren) gets new hour totals.
# The "Multi-Step form" task got broken down into more granular sub-steps, which pushed the
# previously-empty rows 22-25 down into use by the tasks that used to live in rows 18-21.

$ws.Range("B15").Value = "Multi-Step form informatie verzamelen"
$ws.Range("D15").Value = 4

$ws.Range("B16").Value = "Multi-Step form start"
$ws.Range("D16").Value = 5

$ws.Range("B17").Value = "Multi-Step form opmaak stap 1 en stap 2"
$ws.Range("D17").Value = 5

$ws.Range("B18").Value = "Multi-Step form Knoppen functioneel"
$ws.Range("D18").Value = 7

$ws.Range("B19").Value = "Multi-Step form Resultaat geprogrammeerd "
$ws.Range("D19").Value = 5

$ws.Range("B20").Value = "Multi-Step form errors + comments"
$ws.Range("D20").Value = 9

$ws.Range("B21").Value = "Multi-Step form gereed"
$ws.Range("D21").Value = 0.5

$ws.Range("B22").Value = "Systeem Test uitvoeren en bespreken"
$ws.Range("D22").Value = 4

$ws.Range("B23").Value = "Verbetervoorstellen uitwerken"
$ws.Range("D23").Value = 5

$ws.Range("B24").Value = "Alle verbetervoorstellen uitoefenen"
$ws.Range("D24").Value = 10

$ws.Range("B25").Value = "Examenportfolio uitwerken"
$ws.Range("D25").Value = 5.5

# Reflect the author's final selection / scroll position on the Planning sheet.
$ws.Range("D26").Select()
